$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row 21 with "Total" / "Unique" header labels for the train trip stats
$ws.Range("B21").Value = "Total"
$ws.Range("C21").Value = "Unique"

# Add new train-trip rows below the existing bus-stop ridership row
$ws.Range("A23").Value = "R5 Paoli-Thorndale PAO"
$ws.Range("A24").Value = "100 Norristown Speed Line"

# Update the view to match the new selection/scroll position
$excel.ActiveWindow.ScrollRow = 16
$ws.Range("A24").Select()
